$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B format as text so long digit-strings are not coerced to numbers
$ws.Range("B1:B63").NumberFormat = "@"

$ws.Cells.Item(1, 1).Value = 'Фото'
$ws.Cells.Item(1, 2).Value = 'Хэш'
$ws.Cells.Item(1, 3).Value = 'Время обработки'
$ws.Cells.Item(1, 4).Value = 'Хэммингово расстояние'

$ws.Cells.Item(2, 1).Value = 'red_water.jpg'
$ws.Cells.Item(2, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0

$ws.Cells.Item(3, 1).Value = 'red_chb.jpg'
$ws.Cells.Item(3, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(3, 3).Value = 0.015634
$ws.Cells.Item(3, 4).Value = 0

$ws.Cells.Item(4, 1).Value = 'noisy_image_1_110.jpg'
$ws.Cells.Item(4, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0

$ws.Cells.Item(5, 1).Value = 'red.jpg'
$ws.Cells.Item(5, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 0

$ws.Cells.Item(6, 1).Value = 'red_1.jpg'
$ws.Cells.Item(6, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0

$ws.Cells.Item(7, 1).Value = 'red_2.jpg'
$ws.Cells.Item(7, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0

$ws.Cells.Item(8, 1).Value = 'red_3.jpg'
$ws.Cells.Item(8, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0

$ws.Cells.Item(9, 1).Value = 'red_4.jpg'
$ws.Cells.Item(9, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0

$ws.Cells.Item(10, 1).Value = 'red_5.jpg'
$ws.Cells.Item(10, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0

$ws.Cells.Item(11, 1).Value = 'red_6.jpg'
$ws.Cells.Item(11, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0

$ws.Cells.Item(12, 1).Value = 'red_32.jpg'
$ws.Cells.Item(12, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 0

$ws.Cells.Item(13, 1).Value = 'red_33.jpg'
$ws.Cells.Item(13, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 0

$ws.Cells.Item(14, 1).Value = 'red_34.jpg'
$ws.Cells.Item(14, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(14, 3).Value = 0.015618
$ws.Cells.Item(14, 4).Value = 0

$ws.Cells.Item(15, 1).Value = 'red_35.jpg'
$ws.Cells.Item(15, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0

$ws.Cells.Item(16, 1).Value = 'red_36.jpg'
$ws.Cells.Item(16, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 0

$ws.Cells.Item(17, 1).Value = 'red_37.jpg'
$ws.Cells.Item(17, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 0

$ws.Cells.Item(18, 1).Value = 'red_38.jpg'
$ws.Cells.Item(18, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 0

$ws.Cells.Item(19, 1).Value = 'red_39.jpg'
$ws.Cells.Item(19, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 0

$ws.Cells.Item(20, 1).Value = 'red_40.jpg'
$ws.Cells.Item(20, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0

$ws.Cells.Item(21, 1).Value = 'red_41.jpg'
$ws.Cells.Item(21, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = 0

$ws.Cells.Item(22, 1).Value = 'red_42.jpg'
$ws.Cells.Item(22, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 0

$ws.Cells.Item(23, 1).Value = 'red_col_10.jpg'
$ws.Cells.Item(23, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 0

$ws.Cells.Item(24, 1).Value = 'red_col_11.jpg'
$ws.Cells.Item(24, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0

$ws.Cells.Item(25, 1).Value = 'red_col_12.jpg'
$ws.Cells.Item(25, 2).Value = '1110001111000001100000011000000110000001100000011100001111100011'
$ws.Cells.Item(25, 3).Value = 0.015618
$ws.Cells.Item(25, 4).Value = 0

$ws.Cells.Item(26, 1).Value = 'red_col_13.jpg'
$ws.Cells.Item(26, 2).Value = '1110001111000001100000011000000110000001100010011100001111100011'
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 1

$ws.Cells.Item(27, 1).Value = 'red_col_14.jpg'
$ws.Cells.Item(27, 2).Value = '1110001111000001100000011000000110000001100010011100001111100011'
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 1

$ws.Cells.Item(28, 1).Value = 'red_col_15.jpg'
$ws.Cells.Item(28, 2).Value = '1110001111000001100000011000000110000001100010011100001111100011'
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 1

$ws.Cells.Item(29, 1).Value = 'red_col_16.jpg'
$ws.Cells.Item(29, 2).Value = '1110001111000001100000011000000110000001100010011100001111100011'
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 1

$ws.Cells.Item(30, 1).Value = 'red_col_17.jpg'
$ws.Cells.Item(30, 2).Value = '1110001111000001100000011000000110000001100010011100001111100011'
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 1

$ws.Cells.Item(31, 1).Value = 'red_col_18.jpg'
$ws.Cells.Item(31, 2).Value = '1110001111000001100000011000000110000001100010011100001111100011'
$ws.Cells.Item(31, 3).Value = 0.015625
$ws.Cells.Item(31, 4).Value = 1

$ws.Cells.Item(32, 1).Value = 'red_col_19.jpg'
$ws.Cells.Item(32, 2).Value = '1110001111000001100000011000000110000001100010011100001111100011'
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 1

$ws.Cells.Item(33, 1).Value = 'red_col_20.jpg'
$ws.Cells.Item(33, 2).Value = '1110001111000001100000011000000110000001100010011100001111100011'
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 1

$ws.Cells.Item(34, 1).Value = 'new_1.jpg'
$ws.Cells.Item(34, 2).Value = '1111111111100111110001111100001111000111110001111110011111111111'
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 22

$ws.Cells.Item(35, 1).Value = 'new_2.jpg'
$ws.Cells.Item(35, 2).Value = '1110001111001001101000011110000111101101100010011100100111100011'
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 11

$ws.Cells.Item(36, 1).Value = 'new_3.jpg'
$ws.Cells.Item(36, 2).Value = '1111011111000011100000011110000110000001100000011100001111110111'
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 7

$ws.Cells.Item(37, 1).Value = 'new_4.jpg'
$ws.Cells.Item(37, 2).Value = '1111111111010011100100011001000011110000111110011111001111110111'
$ws.Cells.Item(37, 3).Value = 0.015636
$ws.Cells.Item(37, 4).Value = 20

$ws.Cells.Item(38, 1).Value = 'new_5.jpg'
$ws.Cells.Item(38, 2).Value = '1111101110000001101100001011000010011001100110011101101111111111'
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 18

$ws.Cells.Item(39, 1).Value = 'new_6.jpg'
$ws.Cells.Item(39, 2).Value = '1111111111111011100000111000011110000111100001111100111111111111'
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 19

$ws.Cells.Item(40, 1).Value = 'new_7.jpg'
$ws.Cells.Item(40, 2).Value = '1111111111000011100000011000000110000001100000011100001111111111'
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 7

$ws.Cells.Item(41, 1).Value = 'new_8.jpg'
$ws.Cells.Item(41, 2).Value = '1111111111000101100001011000111110001111100011111101111111111111'
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = 20

$ws.Cells.Item(42, 1).Value = 'new_9.jpg'
$ws.Cells.Item(42, 2).Value = '1111111111111111110001111100001110000011110000111110011111111111'
$ws.Cells.Item(42, 3).Value = 0.015631
$ws.Cells.Item(42, 4).Value = 21

$ws.Cells.Item(43, 1).Value = 'new_10.jpg'
$ws.Cells.Item(43, 2).Value = '1111111111000001100000011000100110001101100010011100001111111111'
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(43, 4).Value = 10

$ws.Cells.Item(44, 1).Value = 'new_11.jpg'
$ws.Cells.Item(44, 2).Value = '1111111111000011100000111000001111001011111100111110011111111111'
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 4).Value = 18

$ws.Cells.Item(45, 1).Value = 'new_12.jpg'
$ws.Cells.Item(45, 2).Value = '1111111111000001100000011000010110001101100011011100001111111111'
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = 11

$ws.Cells.Item(46, 1).Value = 'new_13.jpg'
$ws.Cells.Item(46, 2).Value = '1111011111101011110000111000000111001011110010111110101111100111'
$ws.Cells.Item(46, 3).Value = 0.015619
$ws.Cells.Item(46, 4).Value = 16

$ws.Cells.Item(47, 1).Value = 'new_14.jpg'
$ws.Cells.Item(47, 2).Value = '1111111111001011100001111000001110000001100111111111111111111111'
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 19

$ws.Cells.Item(48, 1).Value = 'new_15.jpg'
$ws.Cells.Item(48, 2).Value = '1111011111000011100011011010010100000101100011111100111111101111'
$ws.Cells.Item(48, 3).Value = 0.01563
$ws.Cells.Item(48, 4).Value = 16

$ws.Cells.Item(49, 1).Value = 'new_16.jpg'
$ws.Cells.Item(49, 2).Value = '1110011111000011110010011100100111001001110010011100001111100011'
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 10

$ws.Cells.Item(50, 1).Value = 'new_17.jpg'
$ws.Cells.Item(50, 2).Value = '1111000111111001110000011100010111000001110000011100101111000011'
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 12

$ws.Cells.Item(51, 1).Value = 'new_18.jpg'
$ws.Cells.Item(51, 2).Value = '1111111111001011100000010000010100000001100011111101111111111111'
$ws.Cells.Item(51, 3).Value = 0.017701
$ws.Cells.Item(51, 4).Value = 17

$ws.Cells.Item(52, 1).Value = 'new_19.jpg'
$ws.Cells.Item(52, 2).Value = '1111111111001011100000010000010100000001100011111101111111111111'
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 17

$ws.Cells.Item(53, 1).Value = 'new_20.jpg'
$ws.Cells.Item(53, 2).Value = '1111111111001011100001011000111110001111100011111111111111111111'
$ws.Cells.Item(53, 3).Value = 0.01737
$ws.Cells.Item(53, 4).Value = 22

$ws.Cells.Item(54, 1).Value = 'new_21.jpg'
$ws.Cells.Item(54, 2).Value = '1111101111001001100001010010000000000101000011011100110111101011'
$ws.Cells.Item(54, 3).Value = 0.008822
$ws.Cells.Item(54, 4).Value = 16

$ws.Cells.Item(55, 1).Value = 'new_22.jpg'
$ws.Cells.Item(55, 2).Value = '1111111111110011110000111000000110000001110000111110001111111111'
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 4).Value = 14

$ws.Cells.Item(56, 1).Value = 'new_23.jpg'
$ws.Cells.Item(56, 2).Value = '1101111110011111001100001011000000110000101100001111101111111111'
$ws.Cells.Item(56, 3).Value = 0.015628
$ws.Cells.Item(56, 4).Value = 29

$ws.Cells.Item(57, 1).Value = 'new_24.jpg'
$ws.Cells.Item(57, 2).Value = '1110011111000011110000111100001111000011110000111100011111100111'
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(57, 4).Value = 12

$ws.Cells.Item(58, 1).Value = 'new_25.jpg'
$ws.Cells.Item(58, 2).Value = '1111111111110011100000011000000100000001100010011111101111111111'
$ws.Cells.Item(58, 3).Value = 0.015702
$ws.Cells.Item(58, 4).Value = 14

$ws.Cells.Item(59, 1).Value = 'new_26.jpg'
$ws.Cells.Item(59, 2).Value = '1111111111000111110000111100001111000011110000111100011111101111'
$ws.Cells.Item(59, 3).Value = 0.031236
$ws.Cells.Item(59, 4).Value = 16

$ws.Cells.Item(60, 1).Value = 'new_27.jpg'
$ws.Cells.Item(60, 2).Value = '1111101111000011100001111000000110000001111000111111001111111111'
$ws.Cells.Item(60, 3).Value = 0.031272
$ws.Cells.Item(60, 4).Value = 13

$ws.Cells.Item(61, 1).Value = 'new_28.jpg'
$ws.Cells.Item(61, 2).Value = '1110001111000001100000000000010000000100100011011100110111100011'
$ws.Cells.Item(61, 3).Value = 0.03123
$ws.Cells.Item(61, 4).Value = 12

$ws.Cells.Item(62, 1).Value = 'new_29.jpg'
$ws.Cells.Item(62, 2).Value = '1111101111001011100001111000001100000011100000111000011111111111'
$ws.Cells.Item(62, 3).Value = 0.034884
$ws.Cells.Item(62, 4).Value = 15

$ws.Cells.Item(63, 1).Value = 'new_30.jpg'
$ws.Cells.Item(63, 2).Value = '1111111111110111110000111100001111000011110000111111011111111111'
$ws.Cells.Item(63, 3).Value = 0.037465
$ws.Cells.Item(63, 4).Value = 21

